# Actualización automática del tracker
# Append 5 new result rows (92-96) to the tracker sheet, matching the
# columns: event_id, fecha, jugador_A, jugador_B, pronostico, cuota, resultado, profit
# The two new matches are still pending (resultado/profit), so those are
# left blank, same as other not-yet-settled rows already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TrackerRow($Row, $EventId, $Fecha, $JugadorA, $JugadorB, $Pronostico, $Cuota) {
    $ws.Cells.Item($Row, 1).Value = $EventId

    # Force column B to stay plain text ("2025-10-05") instead of being
    # auto-converted into a date serial number.
    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 2).Value = $Fecha
    $ws.Cells.Item($Row, 2).Style = "Normal"

    $ws.Cells.Item($Row, 3).Value = $JugadorA
    $ws.Cells.Item($Row, 4).Value = $JugadorB
    $ws.Cells.Item($Row, 5).Value = $Pronostico
    $ws.Cells.Item($Row, 6).Value = $Cuota

    # resultado / profit: match still pending, leave as blank text cells
    # (same empty-string convention already used by other pending rows).
    $ws.Cells.Item($Row, 7).Value = "'"
    $ws.Cells.Item($Row, 7).Style = "Normal"

    $ws.Cells.Item($Row, 8).Value = "'"
    $ws.Cells.Item($Row, 8).Style = "Normal"
}

Add-TrackerRow 92 14826685 "2025-10-05" "Anastasia Zakharova" "Bianca Andreescu" "Gana Anastasia Zakharova" 3.25
Add-TrackerRow 93 14828414 "2025-10-05" "Florian Broska" "Dominik Kellovsky" "Gana Dominik Kellovsky" 2.2
Add-TrackerRow 94 14828413 "2025-10-05" "Max Basing" "Filip Jeff Planinsek" "Gana Filip Jeff Planinsek" 2.38
Add-TrackerRow 95 14828405 "2025-10-05" "Dali Blanch" "Adrian Oetzbach" "Gana Adrian Oetzbach" 2.62
Add-TrackerRow 96 14827885 "2025-10-05" "Evan Zhu" "Petr Bar Biryukov" "Gana Evan Zhu" 3
